# Update the "想去人数" (want-to-go count) figures in both the "展览"
# sheet and the aggregated "全部类型" sheet, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 204
$ws1.Range("F13").Value = 789
$ws1.Range("F15").Value = 558
$ws1.Range("F17").Value = 1306
$ws1.Range("F20").Value = 1123
$ws1.Range("F21").Value = 2810
$ws1.Range("F22").Value = 1322
$ws1.Range("F23").Value = 661
$ws1.Range("F27").Value = 978
$ws1.Range("F29").Value = 1582
$ws1.Range("F32").Value = 1345

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F18").Value = 204
$ws4.Range("F26").Value = 789
$ws4.Range("F28").Value = 558
$ws4.Range("F30").Value = 1306
$ws4.Range("F33").Value = 1123
$ws4.Range("F34").Value = 0
$ws4.Range("F35").Value = 1322
$ws4.Range("F36").Value = 661
$ws4.Range("F42").Value = 978
$ws4.Range("F44").Value = 1582
$ws4.Range("F47").Value = 1345
